# Update TPM-derived values on the active worksheet ("update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 21.33926
$ws.Range("N2").Value = 64.01778
$ws.Range("O2").Value = 0.4398914187744692
$ws.Range("P2").Value = 0.4398914187744692
$ws.Range("Q2").Value = 1.08919850892
$ws.Range("R2").Value = 9.802786580280001
$ws.Range("S2").Value = 0.4398914187744692
$ws.Range("T2").Value = 0.4398914187744692

# Row 3
$ws.Range("O3").Value = 0.23906065069302
$ws.Range("P3").Value = 0.23906065069302
$ws.Range("S3").Value = 0.23906065069302
$ws.Range("T3").Value = 0.23906065069302

# Row 4
$ws.Range("M4").Value = 10.59425366666667
$ws.Range("N4").Value = 31.782761
$ws.Range("O4").Value = 0.2183918878295978
$ws.Range("P4").Value = 0.2183918878295978
$ws.Range("Q4").Value = 0.5407518956540001
$ws.Range("R4").Value = 4.866767060886001
$ws.Range("S4").Value = 0.2183918878295978
$ws.Range("T4").Value = 0.2183918878295978

# Row 5
$ws.Range("M5").Value = 4.979874333333333
$ws.Range("N5").Value = 14.939623
$ws.Range("O5").Value = 0.102656042702913
$ws.Range("P5").Value = 0.102656042702913
$ws.Range("Q5").Value = 0.254182745722
$ws.Range("R5").Value = 2.287644711498
$ws.Range("S5").Value = 0.102656042702913
$ws.Range("T5").Value = 0.102656042702913
